$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text. Some of the "price" strings in
# column D look like plain decimal numbers (e.g. "1.00", "598.24") and
# Excel would silently coerce/renormalize them (dropping trailing zeros,
# switching to scientific notation, etc.) if assigned directly. Forcing
# the cell to the "@" (Text) number format before the assignment keeps
# the exact original string, and resetting the style back to "Normal"
# afterwards avoids leaving a stray formatting difference behind.
function Set-TextCell {
    param(
        [string]$CellRef,
        [string]$Val,
        [bool]$ForceText = $false
    )
    $c = $ws.Range($CellRef)
    if ($ForceText) {
        $c.NumberFormat = "@"
        $c.Value = $Val
        $c.Style = "Normal"
    } else {
        $c.Value = $Val
    }
}

# --- Updated price / 1h-volume figures (refreshed crypto snapshot) ---
Set-TextCell "D2"  "67.584.58"  $false
Set-TextCell "E2"  "  -0.27%  " $false
Set-TextCell "D3"  "3.765.63"   $false
Set-TextCell "E3"  "  -0.93%  " $false
Set-TextCell "D4"  "1.00"       $true
Set-TextCell "E4"  "  +0.35%  " $false
Set-TextCell "D5"  "598.24"     $true
Set-TextCell "E5"  "  +0.21%  " $false
Set-TextCell "D6"  "162.75"     $true
Set-TextCell "E6"  "  -2.68%  " $false
Set-TextCell "D7"  "3.761.13"   $false
Set-TextCell "E8"  "  +0.03%  " $false
Set-TextCell "E9"  "  -1.44%  " $false
Set-TextCell "E10" "  -3.02%  " $false
Set-TextCell "D11" "0.444"      $true
Set-TextCell "E11" "  -1.30%  " $false
Set-TextCell "D12" "6.55"       $true
Set-TextCell "E12" "  +3.98%  " $false
Set-TextCell "E13" "  -4.05%  " $false
Set-TextCell "D14" "35.14"      $true
Set-TextCell "E14" "  -2.26%  " $false
Set-TextCell "D15" "4.396.05"   $false
Set-TextCell "D16" "3.778.90"   $false
Set-TextCell "E16" "  -1.14%  " $false
Set-TextCell "D17" "67.648.92"  $false
Set-TextCell "E17" "  -0.24%  " $false
Set-TextCell "D18" "18.25"      $true
Set-TextCell "E18" "  -1.57%  " $false
Set-TextCell "E19" "  +1.63%  " $false
Set-TextCell "D20" "6.98"       $true
Set-TextCell "E20" "  -1.52%  " $false
Set-TextCell "D21" "458.55"     $true
Set-TextCell "E21" "  -0.67%  " $false
Set-TextCell "D22" "9.46"       $true
Set-TextCell "E22" "  -4.53%  " $false
Set-TextCell "D23" "0.691"      $true
Set-TextCell "E23" "  -1.50%  " $false
Set-TextCell "D24" "82.51"      $true
Set-TextCell "E24" "  -0.93%  " $false
Set-TextCell "D25" "0.0000142"  $true
Set-TextCell "E25" "  -6.16%  " $false
Set-TextCell "D26" "11.85"      $true
Set-TextCell "E26" "  -1.97%  " $false
Set-TextCell "E27" "  +0.00%  " $false
Set-TextCell "D28" "2.07"       $true
Set-TextCell "E28" "  -1.49%  " $false
Set-TextCell "D29" "9.81"       $true
Set-TextCell "E29" "  -1.88%  " $false
Set-TextCell "D30" "3.918.74"   $false
Set-TextCell "E30" "  -0.72%  " $false
Set-TextCell "D31" "2.21"       $true
Set-TextCell "E31" "  -1.35%  " $false
Set-TextCell "D32" "7.20"       $true
Set-TextCell "E32" "  -1.82%  " $false
Set-TextCell "D33" "2.57"       $true
Set-TextCell "E33" "  -7.11%  " $false
Set-TextCell "D34" "28.81"      $true
Set-TextCell "E34" "  -2.60%  " $false
Set-TextCell "E35" "  +0.07%  " $false
Set-TextCell "D36" "8.91"       $true
Set-TextCell "E36" "  -1.61%  " $false
Set-TextCell "D37" "0.0986"     $true
Set-TextCell "E37" "  -1.44%  " $false
Set-TextCell "D38" "0.141"      $true
Set-TextCell "E38" "  +2.64%  " $false
Set-TextCell "D39" "5.75"       $true
Set-TextCell "E39" "  -0.46%  " $false
Set-TextCell "D40" "0.974"      $true
Set-TextCell "E40" "  -2.54%  " $false
Set-TextCell "E41" "  -6.01%  " $false
Set-TextCell "D42" "1.00"       $true
Set-TextCell "E42" "  +0.16%  " $false
Set-TextCell "E44" "  -1.84%  " $false
Set-TextCell "D45" "43.02"      $true
Set-TextCell "E45" "  +0.42%  " $false
Set-TextCell "D46" "152.24"     $true
Set-TextCell "E46" "  +3.24%  " $false
Set-TextCell "D47" "0.293"      $true
Set-TextCell "E47" "  -2.90%  " $false
Set-TextCell "D49" "8.25"       $true
Set-TextCell "E49" "  -1.00%  " $false
Set-TextCell "D51" "1.83"       $true
Set-TextCell "E51" "  -0.55%  " $false

# --- Rows 48 & 50 swapped ranking positions: Notcoin <-> ONDO ---
Set-TextCell "B48" "ONDO" $false
Set-TextCell "C48" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo" $false
Set-TextCell "D48" "1.37" $true
Set-TextCell "E48" "  +0.45%  " $false

Set-TextCell "B50" "Notcoin" $false
Set-TextCell "C50" "https://coinranking.com/coin/2L2Y4ghjj+notcoin-not" $false
Set-TextCell "D50" "0.0264" $true
Set-TextCell "E50" "  +79.08%  " $false
